$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.900.42"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.549.42"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.20"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.770.46"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.552.18"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("E14").Value = "  +0.96%  "

$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.906.70"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.64"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.33"
$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.68"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("E33").Value = "  +4.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.411.98"
$ws.Range("E34").Value = "  +2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +2.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("E36").Value = "  -0.66%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("E42").Value = "  +3.72%  "

$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.54"
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.684.81"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.30"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +5.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  +0.46%  "
